$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.964.54"
$ws.Range("E2").Value = "  -2.49%  "
$ws.Range("D3").Value = "3.493.39"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "584.15"
$ws.Range("E5").Value = "  -1.79%  "
$ws.Range("D6").Value = "173.05"
$ws.Range("E6").Value = "  -4.82%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "3.488.96"
$ws.Range("E8").Value = "  +0.38%  "
$ws.Range("B9").Value = "XRP"
$ws.Range("C9").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D9").Value = "0.595"
$ws.Range("E9").Value = "  -4.01%  "
$ws.Range("D10").Value = "0.132"
$ws.Range("E10").Value = "  -6.16%  "
$ws.Range("D11").Value = "6.85"
$ws.Range("E11").Value = "  -2.08%  "
$ws.Range("D12").Value = "0.411"
$ws.Range("E12").Value = "  -4.14%  "
$ws.Range("D13").Value = "4.094.73"
$ws.Range("E13").Value = "  +0.42%  "
$ws.Range("E14").Value = "  +0.35%  "
$ws.Range("D15").Value = "29.99"
$ws.Range("E15").Value = "  -6.56%  "
$ws.Range("D16").Value = "66.097.79"
$ws.Range("E16").Value = "  -2.27%  "
$ws.Range("D17").Value = "0.0000172"
$ws.Range("E17").Value = "  -3.27%  "
$ws.Range("D18").Value = "3.482.00"
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("D19").Value = "5.94"
$ws.Range("E19").Value = "  -4.10%  "
$ws.Range("D20").Value = "13.90"
$ws.Range("E20").Value = "  -1.50%  "
$ws.Range("D21").Value = "367.28"
$ws.Range("E21").Value = "  -6.98%  "
$ws.Range("E22").Value = "  -2.07%  "
$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").Value = "0.541"
$ws.Range("E23").Value = "  +0.37%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D25").Value = "72.09"
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").Value = "9.73"
$ws.Range("E27").Value = "  -6.04%  "
$ws.Range("E28").Value = "  +0.90%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").Value = "24.23"
$ws.Range("E30").Value = "  +2.69%  "
$ws.Range("D31").Value = "5.79"
$ws.Range("E31").Value = "  -5.65%  "
$ws.Range("E32").Value = "  -3.19%  "
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("D34").Value = "1.30"
$ws.Range("E34").Value = "  -7.52%  "
$ws.Range("E35").Value = "  -3.52%  "
$ws.Range("E36").Value = "  -1.04%  "
$ws.Range("D37").Value = "29.69"
$ws.Range("E37").Value = "  +12.96%  "
$ws.Range("D38").Value = "159.14"
$ws.Range("E38").Value = "  -1.65%  "
$ws.Range("D39").Value = "0.888"
$ws.Range("E39").Value = "  -0.52%  "
$ws.Range("E40").Value = "  -4.47%  "
$ws.Range("D41").Value = "2.795.35"
$ws.Range("E41").Value = "  +1.85%  "
$ws.Range("E42").Value = "  -10.63%  "
$ws.Range("E43").Value = "  -5.72%  "
$ws.Range("D44").Value = "6.35"
$ws.Range("E44").Value = "  -5.93%  "
$ws.Range("D45").Value = "0.0689"
$ws.Range("E45").Value = "  -3.89%  "
$ws.Range("D46").Value = "39.92"
$ws.Range("D47").Value = "24.25"
$ws.Range("E47").Value = "  -7.79%  "
$ws.Range("D48").Value = "0.0290"
$ws.Range("E48").Value = "  -3.13%  "
$ws.Range("D49").Value = "307.47"
$ws.Range("E49").Value = "  -6.69%  "
$ws.Range("D50").Value = "0.823"
$ws.Range("E50").Value = "  -2.85%  "
$ws.Range("E51").Value = "  -3.77%  "
